# Update the NATMI ligand-receptor table with re-run (new TPM) results.
# The "Sending cluster" groups change from {FAPs, MuSCs} to {ECs, FAPs}
# (rows 2-4 now report ECs as sender, rows 5-7 report FAPs as sender),
# while the "Target cluster" rotation (ECs / FAPs / MuSCs) and the
# Ligand/Receptor symbols (Epgn / Egfr) stay the same. All of the
# downstream expression / specificity statistics are refreshed to match.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value2 = "ECs"
$ws.Range("B2").Value2 = "Epgn"
$ws.Range("C2").Value2 = "Egfr"
$ws.Range("D2").Value2 = "ECs"
$ws.Range("E2").Value2 = 2
$ws.Range("F2").Value2 = 0.6666666666666666
$ws.Range("G2").Value2 = 0.128109
$ws.Range("H2").Value2 = 0.384327
$ws.Range("I2").Value2 = 0.7686801351245942
$ws.Range("J2").Value2 = 0.7686801351245942
$ws.Range("K2").Value2 = 3
$ws.Range("L2").Value2 = 1
$ws.Range("M2").Value2 = 0.428743
$ws.Range("N2").Value2 = 1.286229
$ws.Range("O2").Value2 = 0.00412050394863168
$ws.Range("P2").Value2 = 0.00412050394863168
$ws.Range("Q2").Value2 = 0.05492583698700001
$ws.Range("R2").Value2 = 0.494332532883
$ws.Range("S2").Value2 = 0.003167349532015624
$ws.Range("T2").Value2 = 0.003167349532015624

$ws.Range("A3").Value2 = "ECs"
$ws.Range("B3").Value2 = "Epgn"
$ws.Range("C3").Value2 = "Egfr"
$ws.Range("D3").Value2 = "FAPs"
$ws.Range("E3").Value2 = 2
$ws.Range("F3").Value2 = 0.6666666666666666
$ws.Range("G3").Value2 = 0.128109
$ws.Range("H3").Value2 = 0.384327
$ws.Range("I3").Value2 = 0.7686801351245942
$ws.Range("J3").Value2 = 0.7686801351245942
$ws.Range("K3").Value2 = 3
$ws.Range("L3").Value2 = 1
$ws.Range("M3").Value2 = 80.22623699999998
$ws.Range("N3").Value2 = 240.678711
$ws.Range("O3").Value2 = 0.7710272268990069
$ws.Range("P3").Value2 = 0.7710272268990069
$ws.Range("Q3").Value2 = 10.277702995833
$ws.Range("R3").Value2 = 92.49932696249698
$ws.Range("S3").Value2 = 0.5926733129574698
$ws.Range("T3").Value2 = 0.5926733129574698

$ws.Range("A4").Value2 = "ECs"
$ws.Range("B4").Value2 = "Epgn"
$ws.Range("C4").Value2 = "Egfr"
$ws.Range("D4").Value2 = "MuSCs"
$ws.Range("E4").Value2 = 2
$ws.Range("F4").Value2 = 0.6666666666666666
$ws.Range("G4").Value2 = 0.128109
$ws.Range("H4").Value2 = 0.384327
$ws.Range("I4").Value2 = 0.7686801351245942
$ws.Range("J4").Value2 = 0.7686801351245942
$ws.Range("K4").Value2 = 3
$ws.Range("L4").Value2 = 1
$ws.Range("M4").Value2 = 23.39612766666667
$ws.Range("N4").Value2 = 70.188383
$ws.Range("O4").Value2 = 0.2248522691523614
$ws.Range("P4").Value2 = 0.2248522691523614
$ws.Range("Q4").Value2 = 2.997254519249
$ws.Range("R4").Value2 = 26.975290673241
$ws.Range("S4").Value2 = 0.1728394726351088
$ws.Range("T4").Value2 = 0.1728394726351088

$ws.Range("A5").Value2 = "FAPs"
$ws.Range("B5").Value2 = "Epgn"
$ws.Range("C5").Value2 = "Egfr"
$ws.Range("D5").Value2 = "ECs"
$ws.Range("E5").Value2 = 1
$ws.Range("F5").Value2 = 0.3333333333333333
$ws.Range("G5").Value2 = 0.038552
$ws.Range("H5").Value2 = 0.115656
$ws.Range("I5").Value2 = 0.2313198648754057
$ws.Range("J5").Value2 = 0.2313198648754058
$ws.Range("K5").Value2 = 3
$ws.Range("L5").Value2 = 1
$ws.Range("M5").Value2 = 0.428743
$ws.Range("N5").Value2 = 1.286229
$ws.Range("O5").Value2 = 0.00412050394863168
$ws.Range("P5").Value2 = 0.00412050394863168
$ws.Range("Q5").Value2 = 0.016528900136
$ws.Range("R5").Value2 = 0.148760101224
$ws.Range("S5").Value2 = 0.0009531544166160561
$ws.Range("T5").Value2 = 0.0009531544166160562

$ws.Range("A6").Value2 = "FAPs"
$ws.Range("B6").Value2 = "Epgn"
$ws.Range("C6").Value2 = "Egfr"
$ws.Range("D6").Value2 = "FAPs"
$ws.Range("E6").Value2 = 1
$ws.Range("F6").Value2 = 0.3333333333333333
$ws.Range("G6").Value2 = 0.038552
$ws.Range("H6").Value2 = 0.115656
$ws.Range("I6").Value2 = 0.2313198648754057
$ws.Range("J6").Value2 = 0.2313198648754058
$ws.Range("K6").Value2 = 3
$ws.Range("L6").Value2 = 1
$ws.Range("M6").Value2 = 80.22623699999998
$ws.Range("N6").Value2 = 240.678711
$ws.Range("O6").Value2 = 0.7710272268990069
$ws.Range("P6").Value2 = 0.7710272268990069
$ws.Range("Q6").Value2 = 3.092881888823999
$ws.Range("R6").Value2 = 27.835936999416
$ws.Range("S6").Value2 = 0.1783539139415371
$ws.Range("T6").Value2 = 0.1783539139415371

$ws.Range("A7").Value2 = "FAPs"
$ws.Range("B7").Value2 = "Epgn"
$ws.Range("C7").Value2 = "Egfr"
$ws.Range("D7").Value2 = "MuSCs"
$ws.Range("E7").Value2 = 1
$ws.Range("F7").Value2 = 0.3333333333333333
$ws.Range("G7").Value2 = 0.038552
$ws.Range("H7").Value2 = 0.115656
$ws.Range("I7").Value2 = 0.2313198648754057
$ws.Range("J7").Value2 = 0.2313198648754058
$ws.Range("K7").Value2 = 3
$ws.Range("L7").Value2 = 1
$ws.Range("M7").Value2 = 23.39612766666667
$ws.Range("N7").Value2 = 70.188383
$ws.Range("O7").Value2 = 0.2248522691523614
$ws.Range("P7").Value2 = 0.2248522691523614
$ws.Range("Q7").Value2 = 0.9019675138053334
$ws.Range("R7").Value2 = 8.117707624248
$ws.Range("S7").Value2 = 0.0520127965172526
$ws.Range("T7").Value2 = 0.05201279651725261
